# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
#
# This script reproduces the commit's re-ordering of several fixture rows
# (130/132/133, 134/136/137, 139/140 and 142/144/145 — the columns B and
# F:AC, i.e. everything except the row index in col A and the Div/Div
# Original Name/Date columns C:E which are shared within each group) plus
# a handful of odds corrections on rows 192-199.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Group 1: rows 130, 132, 133 -- 3-way cyclic rotation of B:AC
#   new130 = old133 ; new132 = old130 ; new133 = old132
# ---------------------------------------------------------------------
$old130 = $ws.Range("B130:AC130").Value()
$old132 = $ws.Range("B132:AC132").Value()
$old133 = $ws.Range("B133:AC133").Value()

$ws.Range("B130:AC130").Value = $old133
$ws.Range("B132:AC132").Value = $old130
$ws.Range("B133:AC133").Value = $old132

# ---------------------------------------------------------------------
# Group 2: rows 134, 136, 137 -- 3-way cyclic rotation of B:AC
#   new134 = old137 ; new136 = old134 ; new137 = old136
# ---------------------------------------------------------------------
$old134 = $ws.Range("B134:AC134").Value()
$old136 = $ws.Range("B136:AC136").Value()
$old137 = $ws.Range("B137:AC137").Value()

$ws.Range("B134:AC134").Value = $old137
$ws.Range("B136:AC136").Value = $old134
$ws.Range("B137:AC137").Value = $old136

# ---------------------------------------------------------------------
# Group 3: rows 139, 140 -- simple swap of B:AC
# ---------------------------------------------------------------------
$old139 = $ws.Range("B139:AC139").Value()
$old140 = $ws.Range("B140:AC140").Value()

$ws.Range("B139:AC139").Value = $old140
$ws.Range("B140:AC140").Value = $old139

# ---------------------------------------------------------------------
# Group 4: rows 142, 144, 145 -- 3-way cyclic rotation of B:AC
#   new142 = old144 ; new144 = old145 ; new145 = old142
# ---------------------------------------------------------------------
$old142 = $ws.Range("B142:AC142").Value()
$old144 = $ws.Range("B144:AC144").Value()
$old145 = $ws.Range("B145:AC145").Value()

$ws.Range("B142:AC142").Value = $old144
$ws.Range("B144:AC144").Value = $old145
$ws.Range("B145:AC145").Value = $old142

# ---------------------------------------------------------------------
# Odds corrections on rows 192-199 (individual cell updates)
# ---------------------------------------------------------------------

# Row 192
$ws.Range("R192").Value = 1.85
$ws.Range("S192").Value = 1.95
$ws.Range("T192").Value = 2.25
$ws.Range("U192").Value = 1.8
$ws.Range("V192").Value = 2

# Row 193
$ws.Range("U193").Value = 1.825
$ws.Range("V193").Value = 1.975

# Row 194
$ws.Range("R194").Value = 1.875
$ws.Range("S194").Value = 1.925
$ws.Range("U194").Value = 1.925
$ws.Range("V194").Value = 1.875

# Row 195
$ws.Range("N195").Value = 1.65
$ws.Range("O195").Value = 3.6
$ws.Range("P195").Value = 5.25
$ws.Range("R195").Value = 1.825
$ws.Range("S195").Value = 1.975

# Row 196
$ws.Range("R196").Value = 1.975
$ws.Range("S196").Value = 1.825
$ws.Range("U196").Value = 1.925
$ws.Range("V196").Value = 1.875

# Row 197
$ws.Range("N197").Value = 2.875
$ws.Range("O197").Value = 3.25
$ws.Range("P197").Value = 2.45
$ws.Range("Q197").Value = 0
$ws.Range("R197").Value = 2.05
$ws.Range("S197").Value = 1.75

# Row 199
$ws.Range("N199").Value = 2.875
$ws.Range("P199").Value = 2.5
$ws.Range("R199").Value = 2
$ws.Range("S199").Value = 1.8
$ws.Range("U199").Value = 1.85
$ws.Range("V199").Value = 1.95
